$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 0.03819444444444445
$ws.Range("D2").Value = 0.06388888888888888
$ws.Range("E2").Value = "C24"
$ws.Range("F2").Value = 2
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 0.02916666666666667
$ws.Range("D3").Value = 0.05069444444444444
$ws.Range("E3").Value = "D38"
$ws.Range("F3").Value = 1
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 0.01111111111111111
$ws.Range("D4").Value = 0.03263888888888889
$ws.Range("E4").Value = "C12"
$ws.Range("F4").Value = 1
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 0.02430555555555556
$ws.Range("D5").Value = 0.05
$ws.Range("E5").Value = "C13"
$ws.Range("F5").Value = 2
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = 0.0006944444444444445
$ws.Range("D6").Value = 0.02361111111111111
$ws.Range("E6").Value = "D44"
$ws.Range("F6").Value = 3
$ws.Range("B7").Value = 6
$ws.Range("C7").Value = 0.03819444444444445
$ws.Range("D7").Value = 0.06458333333333334
$ws.Range("E7").Value = "E21"
$ws.Range("F7").Value = 1
$ws.Range("B8").Value = 7
$ws.Range("C8").Value = 0.03472222222222222
$ws.Range("D8").Value = 0.05902777777777778
$ws.Range("E8").Value = "A5"
$ws.Range("F8").Value = 2
$ws.Range("B9").Value = 8
$ws.Range("C9").Value = 0.03055555555555555
$ws.Range("D9").Value = 0.05555555555555555
$ws.Range("E9").Value = "B4"
$ws.Range("F9").Value = 1
$ws.Range("B10").Value = 9
$ws.Range("C10").Value = 0.01805555555555555
$ws.Range("D10").Value = 0.04375
$ws.Range("E10").Value = "D34"
$ws.Range("F10").Value = 2
$ws.Range("B11").Value = 10
$ws.Range("C11").Value = 0.009027777777777777
$ws.Range("D11").Value = 0.03194444444444444
$ws.Range("E11").Value = "A12"
$ws.Range("F11").Value = 2
$ws.Range("B12").Value = 11
$ws.Range("C12").Value = 0.004861111111111111
$ws.Range("D12").Value = 0.02847222222222222
$ws.Range("E12").Value = "E1"
$ws.Range("F12").Value = 3
$ws.Range("B13").Value = 12
$ws.Range("C13").Value = 0.02361111111111111
$ws.Range("D13").Value = 0.04861111111111111
$ws.Range("E13").Value = "F60"
$ws.Range("F13").Value = 3
$ws.Range("B14").Value = 13
$ws.Range("C14").Value = 0.04097222222222222
$ws.Range("D14").Value = 0.06597222222222222
$ws.Range("E14").Value = "C3"
$ws.Range("F14").Value = 1
$ws.Range("B15").Value = 14
$ws.Range("C15").Value = 0.00625
$ws.Range("D15").Value = 0.02986111111111111
$ws.Range("E15").Value = "F36"
$ws.Range("F15").Value = 3
$ws.Range("B16").Value = 15
$ws.Range("C16").Value = 0.02291666666666667
$ws.Range("D16").Value = 0.04861111111111111
$ws.Range("E16").Value = "D35"
$ws.Range("F16").Value = 2
$ws.Range("B17").Value = 16
$ws.Range("C17").Value = 0.01805555555555555
$ws.Range("D17").Value = 0.04097222222222222
$ws.Range("E17").Value = "E27"
$ws.Range("F17").Value = 1
$ws.Range("B18").Value = 17
$ws.Range("C18").Value = 0.008333333333333333
$ws.Range("D18").Value = 0.02986111111111111
$ws.Range("E18").Value = "D33"
$ws.Range("F18").Value = 2
$ws.Range("B19").Value = 18
$ws.Range("C19").Value = 0.03402777777777777
$ws.Range("D19").Value = 0.05972222222222223
$ws.Range("E19").Value = "B8"
$ws.Range("F19").Value = 1
$ws.Range("B20").Value = 19
$ws.Range("C20").Value = 0.03055555555555555
$ws.Range("D20").Value = 0.05763888888888889
$ws.Range("E20").Value = "E10"
$ws.Range("F20").Value = 1
$ws.Range("B21").Value = 20
$ws.Range("C21").Value = 0.02291666666666667
$ws.Range("D21").Value = 0.04652777777777778
$ws.Range("E21").Value = "A19"
$ws.Range("F21").Value = 2
